$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 23.86000000000029
$ws.Range("H2").Value = [double]"5.844220662964972e-10"
$ws.Range("I2").Value = [double]"5.844220662964972e-10"
$ws.Range("L2").Value = 52.44626921797015
$ws.Range("M2").Value = '[38.54292359529603, 66.34961484064428]'
$ws.Range("N2").Value = [double]"1.336083688130429e-09"
$ws.Range("O2").Value = [double]"1.336083688130429e-09"
$ws.Range("P2").Value = 1.62897396852804
$ws.Range("Q2").Value = '[1.3270791789938858, 1.9308687580621946]'
$ws.Range("R2").Value = [double]"3.597122599785507e-14"
$ws.Range("S2").Value = [double]"3.597122599785507e-14"
$ws.Range("T2").Value = 58.31606853415646
$ws.Range("U2").Value = '[49.46853113857937, 67.16360592973354]'
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 17.67407407407429
$ws.Range("Y2").Value = 16.52764764764785
$ws.Range("Z2").Value = 18.82050050050073
$ws.Range("F3").Value = 23.86000000000029
$ws.Range("H3").Value = [double]"5.612552522737957e-08"
$ws.Range("I3").Value = [double]"5.612552522737957e-08"
$ws.Range("L3").Value = 40.22008690154149
$ws.Range("M3").Value = '[26.69655432647336, 53.74361947660961]'
$ws.Range("N3").Value = [double]"3.229756400102701e-07"
$ws.Range("O3").Value = [double]"3.229756400102701e-07"
$ws.Range("P3").Value = 1.679289766783733
$ws.Range("Q3").Value = '[1.2767633807381946, 2.081816152829272]'
$ws.Range("R3").Value = [double]"9.009881729582503e-11"
$ws.Range("S3").Value = [double]"9.009881729582503e-11"
$ws.Range("T3").Value = 55.80642211393003
$ws.Range("U3").Value = '[47.703042851118916, 63.909801376741136]'
$ws.Range("X3").Value = 17.48300300300322
$ws.Range("Y3").Value = 15.95443443443463
$ws.Range("Z3").Value = 19.0115715715718
$ws.Range("F4").Value = 23.86000000000029
$ws.Range("H4").Value = [double]"5.270054059813845e-08"
$ws.Range("I4").Value = [double]"5.270054059813845e-08"
$ws.Range("L4").Value = 40.61168978204653
$ws.Range("M4").Value = '[27.98138371155534, 53.24199585253773]'
$ws.Range("N4").Value = [double]"6.116834083691458e-08"
$ws.Range("O4").Value = [double]"6.116834083691458e-08"
$ws.Range("P4").Value = 1.578658170272348
$ws.Range("Q4").Value = '[1.2138686329185786, 1.9434477076261176]'
$ws.Range("R4").Value = [double]"3.202371701149787e-11"
$ws.Range("S4").Value = [double]"3.202371701149787e-11"
$ws.Range("T4").Value = 58.40609722366631
$ws.Range("U4").Value = '[50.23413064972753, 66.5780637976051]'
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 17.86514514514536
$ws.Range("Y4").Value = 16.47987987988008
$ws.Range("Z4").Value = 19.25041041041064
$ws.Range("F5").Value = 23.86000000000029
$ws.Range("H5").Value = [double]"8.433723652778724e-08"
$ws.Range("I5").Value = [double]"8.433723652778724e-08"
$ws.Range("L5").Value = 45.10666840640788
$ws.Range("M5").Value = '[28.933000409898398, 61.28033640291736]'
$ws.Range("N5").Value = [double]"1.151588829406691e-06"
$ws.Range("O5").Value = [double]"1.151588829406691e-06"
$ws.Range("P5").Value = 1.905710858934349
$ws.Range("Q5").Value = '[1.4906055233248878, 2.320816194543811]'
$ws.Range("R5").Value = [double]"5.703437722104354e-12"
$ws.Range("S5").Value = [double]"5.703437722104354e-12"
$ws.Range("T5").Value = 63.42911593227548
$ws.Range("U5").Value = '[54.161578206926166, 72.6966536576248]'
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 16.62318318318339
$ws.Range("Y5").Value = 15.04684684684703
$ws.Range("Z5").Value = 18.19951951951974
$ws.Range("F6").Value = 23.86000000000029
$ws.Range("H6").Value = [double]"2.517030617266158e-08"
$ws.Range("I6").Value = [double]"2.517030617266158e-08"
$ws.Range("L6").Value = 46.75784142324967
$ws.Range("M6").Value = '[31.257334145932333, 62.25834870056701]'
$ws.Range("N6").Value = [double]"2.41100214681822e-07"
$ws.Range("O6").Value = [double]"2.41100214681822e-07"
$ws.Range("P6").Value = 1.79250031285904
$ws.Range("Q6").Value = '[1.42771077550527, 2.157289850212811]'
$ws.Range("R6").Value = [double]"7.187583861423263e-13"
$ws.Range("S6").Value = [double]"7.187583861423263e-13"
$ws.Range("T6").Value = 54.62968441765515
$ws.Range("U6").Value = '[45.51108354360374, 63.74828529170656]'
$ws.Range("V6").Value = [double]"1.110223024625157e-15"
$ws.Range("W6").Value = [double]"1.110223024625157e-15"
$ws.Range("X6").Value = 17.0530930930933
$ws.Range("Y6").Value = 15.66782782782802
$ws.Range("Z6").Value = 18.43835835835859
$ws.Range("F7").Value = 23.86000000000029
$ws.Range("H7").Value = [double]"6.461085000353251e-10"
$ws.Range("I7").Value = [double]"6.461085000353251e-10"
$ws.Range("L7").Value = 49.43952850318682
$ws.Range("M7").Value = '[34.70331752343591, 64.17573948293773]'
$ws.Range("N7").Value = [double]"2.33638988156315e-08"
$ws.Range("O7").Value = [double]"2.33638988156315e-08"
$ws.Range("P7").Value = 1.817658211986887
$ws.Range("Q7").Value = '[1.490605523324887, 2.144710900648888]'
$ws.Range("R7").Value = [double]"1.354472090042691e-14"
$ws.Range("S7").Value = [double]"1.354472090042691e-14"
$ws.Range("T7").Value = 59.09832566447921
$ws.Range("U7").Value = '[50.725898808663075, 67.47075252029535]'
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 16.95755755755776
$ws.Range("Y7").Value = 15.71559559559579
$ws.Range("Z7").Value = 18.19951951951974
$ws.Range("F8").Value = 23.86000000000029
$ws.Range("H8").Value = [double]"4.203436154703866e-09"
$ws.Range("I8").Value = [double]"4.203436154703866e-09"
$ws.Range("L8").Value = 50.11166655066951
$ws.Range("M8").Value = '[33.302622751873585, 66.92071034946544]'
$ws.Range("N8").Value = [double]"3.074612346942729e-07"
$ws.Range("O8").Value = [double]"3.074612346942729e-07"
$ws.Range("P8").Value = 2.006342455445734
$ws.Range("Q8").Value = '[1.6415529180919632, 2.371131992799504]'
$ws.Range("R8").Value = [double]"1.909583602355269e-14"
$ws.Range("S8").Value = [double]"1.909583602355269e-14"
$ws.Range("T8").Value = 53.90272490940151
$ws.Range("U8").Value = '[44.79935651831279, 63.006093300490235]'
$ws.Range("V8").Value = [double]"1.554312234475219e-15"
$ws.Range("W8").Value = [double]"1.554312234475219e-15"
$ws.Range("X8").Value = 16.24104104104124
$ws.Range("Y8").Value = 14.85577577577595
$ws.Range("Z8").Value = 17.62630630630652
$ws.Range("F9").Value = 23.86000000000029
$ws.Range("H9").Value = [double]"3.61801200099432e-06"
$ws.Range("I9").Value = [double]"3.61801200099432e-06"
$ws.Range("L9").Value = 39.32871950985429
$ws.Range("M9").Value = '[22.241929287012255, 56.41550973269632]'
$ws.Range("N9").Value = [double]"3.0593261736378e-05"
$ws.Range("O9").Value = [double]"3.0593261736378e-05"
$ws.Range("P9").Value = 1.855395060678656
$ws.Range("Q9").Value = '[1.352237078121731, 2.35855304323558]'
$ws.Range("R9").Value = [double]"2.3813579996812e-09"
$ws.Range("S9").Value = [double]"2.3813579996812e-09"
$ws.Range("T9").Value = 58.01121042098661
$ws.Range("U9").Value = '[48.36745574251387, 67.65496509945936]'
$ws.Range("V9").Value = [double]"8.881784197001252e-16"
$ws.Range("W9").Value = [double]"8.881784197001252e-16"
$ws.Range("X9").Value = 16.81425425425446
$ws.Range("Y9").Value = 14.90354354354373
$ws.Range("Z9").Value = 18.7249649649652
$ws.Range("F10").Value = 23.86000000000029
$ws.Range("H10").Value = [double]"3.333083808954029e-10"
$ws.Range("I10").Value = [double]"3.333083808954029e-10"
$ws.Range("L10").Value = 52.89864499582566
$ws.Range("M10").Value = '[37.53595359609608, 68.26133639555525]'
$ws.Range("N10").Value = [double]"1.271421989912369e-08"
$ws.Range("O10").Value = [double]"1.271421989912369e-08"
$ws.Range("P10").Value = 1.855395060678656
$ws.Range("Q10").Value = '[1.5283423720166542, 2.182447749340657]'
$ws.Range("R10").Value = [double]"6.661338147750939e-15"
$ws.Range("S10").Value = [double]"6.661338147750939e-15"
$ws.Range("T10").Value = 60.9686511841302
$ws.Range("U10").Value = '[52.227307776713, 69.7099945915474]'
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0
$ws.Range("X10").Value = 16.81425425425446
$ws.Range("Y10").Value = 15.57229229229248
$ws.Range("Z10").Value = 18.05621621621644
$ws.Range("F11").Value = 23.24000000000019
$ws.Range("H11").Value = [double]"4.103928308296645e-11"
$ws.Range("I11").Value = [double]"4.103928308296645e-11"
$ws.Range("L11").Value = 54.21010209746521
$ws.Range("M11").Value = '[40.07624142449497, 68.34396277043545]'
$ws.Range("N11").Value = [double]"8.687848218613681e-10"
$ws.Range("O11").Value = [double]"8.687848218613681e-10"
$ws.Range("P11").Value = 1.389973926813502
$ws.Range("Q11").Value = '[1.1006580868432705, 1.6792897667837332]'
$ws.Range("R11").Value = [double]"1.442179708988078e-12"
$ws.Range("S11").Value = [double]"1.442179708988078e-12"
$ws.Range("T11").Value = 58.98065450734864
$ws.Range("U11").Value = '[50.747781012552664, 67.21352800214461]'
$ws.Range("V11").Value = 0
$ws.Range("W11").Value = 0
$ws.Range("X11").Value = 18.09881881881897
$ws.Range("Y11").Value = 17.02870870870885
$ws.Range("Z11").Value = 19.16892892892909
$ws.Range("F12").Value = 23.24000000000019
$ws.Range("H12").Value = [double]"3.819359939427613e-09"
$ws.Range("I12").Value = [double]"3.819359939427613e-09"
$ws.Range("L12").Value = 47.69800672592151
$ws.Range("M12").Value = '[32.25030591710904, 63.14570753473398]'
$ws.Range("N12").Value = [double]"1.476176711978638e-07"
$ws.Range("O12").Value = [double]"1.476176711978638e-07"
$ws.Range("P12").Value = 1.968605606753964
$ws.Range("Q12").Value = '[1.6163950189641163, 2.320816194543811]'
$ws.Range("R12").Value = [double]"1.13242748511766e-14"
$ws.Range("S12").Value = [double]"1.13242748511766e-14"
$ws.Range("T12").Value = 51.41202158344275
$ws.Range("U12").Value = '[42.83727845285039, 59.98676471403511]'
$ws.Range("V12").Value = [double]"8.881784197001252e-16"
$ws.Range("W12").Value = [double]"8.881784197001252e-16"
$ws.Range("X12").Value = 15.95859859859873
$ws.Range("Y12").Value = 14.65585585585598
$ws.Range("Z12").Value = 17.26134134134149
$ws.Range("F13").Value = 23.24000000000019
$ws.Range("H13").Value = [double]"8.30995816603064e-09"
$ws.Range("I13").Value = [double]"8.30995816603064e-09"
$ws.Range("L13").Value = 45.86239313492663
$ws.Range("M13").Value = '[30.5902799861953, 61.134506283657956]'
$ws.Range("N13").Value = [double]"2.64634216540216e-07"
$ws.Range("O13").Value = [double]"2.64634216540216e-07"
$ws.Range("P13").Value = 1.867974010242579
$ws.Range("Q13").Value = '[1.515763422452732, 2.220184598032427]'
$ws.Range("R13").Value = [double]"6.328271240363392e-14"
$ws.Range("S13").Value = [double]"6.328271240363392e-14"
$ws.Range("T13").Value = 48.25856607768148
$ws.Range("U13").Value = '[39.77927516833208, 56.73785698703088]'
$ws.Range("V13").Value = [double]"6.217248937900877e-15"
$ws.Range("W13").Value = [double]"6.217248937900877e-15"
$ws.Range("X13").Value = 16.33081081081095
$ws.Range("Y13").Value = 15.02806806806819
$ws.Range("Z13").Value = 17.6335535535537
$ws.Range("F14").Value = 23.24000000000019
$ws.Range("H14").Value = [double]"1.486405443174021e-10"
$ws.Range("I14").Value = [double]"1.486405443174021e-10"
$ws.Range("L14").Value = 48.56401675339361
$ws.Range("M14").Value = '[34.06355920997209, 63.06447429681513]'
$ws.Range("N14").Value = [double]"2.432213519476534e-08"
$ws.Range("O14").Value = [double]"2.432213519476534e-08"
$ws.Range("P14").Value = 2.094395102393195
$ws.Range("Q14").Value = '[1.7799213632951174, 2.4088688414912722]'
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 0
$ws.Range("T14").Value = 55.15004938964165
$ws.Range("U14").Value = '[47.36986657438567, 62.930232204897635]'
$ws.Range("V14").Value = 0
$ws.Range("W14").Value = 0
$ws.Range("X14").Value = 16.24104104104124
$ws.Range("Y14").Value = 14.33017017017029
$ws.Range("Z14").Value = 16.65649649649664
